$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new "branchAndBound" worksheet right after "rootnode".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "branchAndBound"

# Header row (B1:E1 first, then A1 - matches the authoring order of the
# shared-string table in the target workbook).
$ws2.Range("B1").Value = "LB"
$ws2.Range("C1").Value = "UB"
$ws2.Range("D1").Value = "active nodes"
$ws2.Range("E1").Value = "gap "
$ws2.Range("A1").Value = "Number of nodes"

# ---------------------------------------------------------------------------
# 2. Update the Benders-iteration log values in column F of "rootnode"
#    (Cross Decomposition iteration summary for the final/best iteration).
# ---------------------------------------------------------------------------
$ws1.Range("F4").Value = "Iteration Value : -1335.7347667524655"
$ws1.Range("F5").Value = "Iteration Bound : -1354.2667698798532"
$ws1.Range("F6").Value = "Objective Value : -1335.736749892465"
$ws1.Range("F8").Value = "Gap : 1.37 %"

# ---------------------------------------------------------------------------
# 3. Finish the branchAndBound header row + data.
# ---------------------------------------------------------------------------
$ws2.Range("F1").Value = "Benders time"
$ws2.Range("G1").Value = "Lagrangean time"

# Row 2 only carries the percentage number format on E2 (no value yet).
$ws2.Range("E2").NumberFormat = "0.00%"

# Data rows.
$ws2.Range("A3").Value = 21
$ws2.Range("B3").Value = -1355.1
$ws2.Range("C3").Value = -1338.88
$ws2.Range("D3").Value = 4
$ws2.Range("E3").NumberFormat = "0.00%"
$ws2.Range("E3").Formula = "=-(C3-B3)/B3"
$ws2.Range("F3").Value = 674
$ws2.Range("G3").Value = 1160

$ws2.Range("A4").Value = 101
$ws2.Range("B4").Value = -1343.34
$ws2.Range("C4").Value = -1338.88
$ws2.Range("D4").Value = 11
$ws2.Range("E4").NumberFormat = "0.00%"
$ws2.Range("E4").Formula = "=-(C4-B4)/B4"
$ws2.Range("F4").Value = 1368
$ws2.Range("G4").Value = 4121

$ws2.Range("A5").Value = 501
$ws2.Range("B5").Value = -1340.32
$ws2.Range("C5").Value = -1338.88
$ws2.Range("D5").Value = 48
$ws2.Range("E5").NumberFormat = "0.00%"
$ws2.Range("E5").Formula = "=-(C5-B5)/B5"
$ws2.Range("F5").Value = 7631
$ws2.Range("G5").Value = 11455

# Column widths (best-fit look-alike; engine rounds to its own pixel grid).
$ws2.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 10.5
$ws2.Columns.Item(6).ColumnWidth = 11.166666666666666
$ws2.Columns.Item(7).ColumnWidth = 14

# ---------------------------------------------------------------------------
# 4. View state: rootnode's old selection moves, and branchAndBound becomes
#    the active / selected sheet + tab.
# ---------------------------------------------------------------------------
$ws1.Range("F5").Select()
$ws2.Range("G9").Select()
$ws2.Activate()
$excel.ActiveWindow.Zoom = 130
